$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Replace the text of the final paragraph (currently the italic
#    "Please create a feature image..." image-prompt) with the
#    "meta description" sentence.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Find.Execute(
    "Please create a feature image that complements the retro theme of the " + [char]34 + "Fruit Super Nova" + [char]34 + " game. The image should showcase a happy Maya warrior wearing glasses. The image should be in cartoon style, with a colorful and eye-catching design. The warrior can be seen holding a basket of brightly colored fruits, with a flaming star (the scatter symbol) just above the basket. The background can feature a starry night sky, with the silhouette of a tropical forest in the distance. The image should convey a fun and exciting atmosphere, inviting players to try out the game and discover its simplicity and entertainment value.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Fruit Super Nova, a simple and fun slot game with fruit-associated jackpots. Play for free and win big prizes.",
    2)

# ------------------------------------------------------------------
# 2. Insert a new paragraph right before that last paragraph holding a
#    bold "Play Fruit Super Nova - Free Online Slot Game" run (plus
#    the leading empty run the rest of the document uses).  Using
#    InsertXML lets us control the run/paragraph structure precisely.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$xmlFrag = '<w:p xmlns:w="' + $wNs + '"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Super Nova - Free Online Slot Game</w:t></w:r></w:p><w:p xmlns:w="' + $wNs + '"></w:p>'
$insertPoint.InsertXML($xmlFrag)

# InsertXML above splits the paragraph, but leaves one spare empty
# paragraph behind it; remove that spare paragraph mark so the new
# bold paragraph sits directly before the (now meta-description) one.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$spacerPara = $lastPara.Previous()
$spacerPara.Range.Delete()

# ------------------------------------------------------------------
# 3. Delete the original "Meta description" paragraph that followed
#    the "Play Fruit Super Nova - Free Online Slot Game" heading.
#    Located by its leading bold "Meta description" text rather than a
#    hard-coded index, in case the paragraph numbering ever shifts.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith("Meta description")) {
        $candidate.Range.Delete()
        break
    }
}
